$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two rows at row 24 to make room, shifting old rows 24-25 to 26-27
$ws.Rows("24:25").Insert()

# New rows 18-22: Students tests.
# Values are entered in this specific order so the shared-string table
# ends up built in the same sequence as the target workbook.
$ws.Range("A18").Value = "Переход к Students"
$ws.Range("A19").Value = "Создание Students"
$ws.Range("A20").Value = "Удаление Students"
$ws.Range("A22").Value = "Поиск Students"
$ws.Range("A21").Value = "Редактирование Students"

$ws.Range("B22").Value = "Can_Search_Students"
$ws.Range("B18").Value = "Can_Navigate_To_Students"
$ws.Range("B19").Value = "Can_Create_Students"
$ws.Range("B20").Value = "Can_Delete_Students"
$ws.Range("B21").Value = "Can_Edit_Student"

$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 0

# New row 24: Student Body Statistics navigation test
$ws.Range("A24").Value = "Переход к Student Body Statistics"
$ws.Range("B24").Value = "Can_Navigate_To_About"
$ws.Range("C24").Value = 0

# Update selection to match the target state
$ws.Range("C24").Select()
